# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-04 10:17:38
#
# Reorders the "Recorded By" email lists (column G) on several session rows,
# updates two fraction-style "Students" counts (H27, H55), and refreshes the
# computed percentage summaries (L10, S15, S16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText {
    # Plain text (email lists, "N/NNN" fraction labels, etc.) - assigning
    # directly through .Value is safe for these and keeps the cell's
    # existing style untouched.
    param($Worksheet, $Row, $Column, $Text)
    $Worksheet.Cells.Item($Row, $Column).Value = $Text
}

function Set-PercentText {
    # Strings like "43.7%" read as a number to Excel's normal Value setter,
    # which would silently convert them to a numeric percentage and
    # (re)stamp a new number-format style on the cell. Routing the literal
    # through a self-quoting formula keeps the cell's existing style intact
    # while still landing the exact text as the cell's value.
    param($Worksheet, $Row, $Column, $Text)
    $Worksheet.Cells.Item($Row, $Column).Formula = '="' + $Text.Replace('"', '""') + '"'
}

# Column G (Recorded By) reorderings
Set-PlainText $ws 2 7 "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
Set-PlainText $ws 3 7 "mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 4 7 "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 5 7 "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
Set-PlainText $ws 6 7 "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 7 7 "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Average Attendance % summary value
Set-PercentText $ws 10 12 "43.7%"

Set-PlainText $ws 12 7 "System, salma.elgendy.std@med.asu.edu.eg"
Set-PlainText $ws 13 7 "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
Set-PlainText $ws 14 7 "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Class Statistics table percentages
Set-PercentText $ws 15 19 "49.2%"
Set-PercentText $ws 16 19 "38.2%"

Set-PlainText $ws 20 7 "user@user.com, nourhan.mostafa@med.asu.edu.eg"

Set-PlainText $ws 24 7 "Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
Set-PlainText $ws 25 7 "youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg"

Set-PlainText $ws 27 7 "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
Set-PlainText $ws 27 8 "20/221"

Set-PlainText $ws 30 7 "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 31 7 "mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 32 7 "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 33 7 "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
Set-PlainText $ws 34 7 "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
Set-PlainText $ws 35 7 "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

Set-PlainText $ws 40 7 "System, salma.elgendy.std@med.asu.edu.eg"
Set-PlainText $ws 41 7 "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
Set-PlainText $ws 42 7 "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

Set-PlainText $ws 48 7 "user@user.com, nourhan.mostafa@med.asu.edu.eg"

Set-PlainText $ws 52 7 "Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
Set-PlainText $ws 53 7 "youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg"

Set-PlainText $ws 55 7 "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
Set-PlainText $ws 55 8 "67/246"
